# Update "想去人数" (interest count) figures in column F across the
# affected sheets, matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 853
$ws.Range("F4").Value = 1071
$ws.Range("F5").Value = 497
$ws.Range("F6").Value = 209
$ws.Range("F7").Value = 649
$ws.Range("F8").Value = 225
$ws.Range("F13").Value = 1740
$ws.Range("F14").Value = 417
$ws.Range("F16").Value = 484
$ws.Range("F17").Value = 194
$ws.Range("F18").Value = 404
$ws.Range("F21").Value = 649
$ws.Range("F24").Value = 949
$ws.Range("F26").Value = 1508
$ws.Range("F27").Value = 254
$ws.Range("F28").Value = 27

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 66

# Sheet "全部类型" (All types - combined list)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 853
$ws.Range("F5").Value = 1071
$ws.Range("F8").Value = 497
$ws.Range("F9").Value = 209
$ws.Range("F10").Value = 649
$ws.Range("F12").Value = 225
$ws.Range("F17").Value = 1740
$ws.Range("F19").Value = 417
$ws.Range("F21").Value = 484
$ws.Range("F22").Value = 194
$ws.Range("F23").Value = 404
$ws.Range("F28").Value = 66
$ws.Range("F29").Value = 649
$ws.Range("F36").Value = 949
$ws.Range("F38").Value = 1508
$ws.Range("F39").Value = 254
$ws.Range("F40").Value = 27
